$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3351.4
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3351.4
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 10054.2
$ws.Range("N17").Value = -10390.2

$ws.Range("H31").Value = 999.5
$ws.Range("I31").Value = 999.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2998.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2768.5

$ws.Range("H61").Value = 483.33334
$ws.Range("I61").Value = 483.33334
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1450.00002
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1278.00002

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1665
$ws.Range("I32").Value = 1665
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1665
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1378

$ws.Range("H61").Value = 8029.7
$ws.Range("I61").Value = 5162.125
$ws.Range("J61").Value = 19500
$ws.Range("K61").Value = 5162.125
$ws.Range("L61").Value = 19500
$ws.Range("M61").Value = -4950.125
$ws.Range("N61").Value = -19924

$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -21400

$ws.Range("H136").Value = 8029.7
$ws.Range("I136").Value = 5162.125
$ws.Range("J136").Value = 19500
$ws.Range("K136").Value = 15486.375
$ws.Range("L136").Value = 58500
$ws.Range("M136").Value = -12936.375
$ws.Range("N136").Value = -63600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4100
$ws.Range("I99").Value = 4100
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4100
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -2602

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H134").Value = 9986.888999999999
$ws.Range("I134").Value = 6346.4
$ws.Range("J134").Value = 14537.5
$ws.Range("K134").Value = 19039.2
$ws.Range("L134").Value = 43612.5
$ws.Range("M134").Value = -16504.2
$ws.Range("N134").Value = -48682.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 1000
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 1000
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -1280

$ws.Range("H80").Value = 25000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 25000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246

$ws.Range("H83").Value = 25000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 25000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232

$ws.Range("H107").Value = 661.38464
$ws.Range("I107").Value = 681.7273
$ws.Range("J107").Value = 549.5
$ws.Range("K107").Value = 681.7273
$ws.Range("L107").Value = 549.5
$ws.Range("M107").Value = 1238.2727
$ws.Range("N107").Value = -4389.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 100
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -524

$ws.Range("H108").Value = 527
$ws.Range("I108").Value = 527
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 1581
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = 1299
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4000
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -3471

$ws.Range("H132").Value = 234923.33
$ws.Range("I132").Value = 297330
$ws.Range("J132").Value = 16500
$ws.Range("K132").Value = 891990
$ws.Range("L132").Value = 49500
$ws.Range("M132").Value = -889460
$ws.Range("N132").Value = -54560

$ws.Range("H141").Value = 112436
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 112436
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 112436
$ws.Range("N141").Value = -122796

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 5000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5226

$ws.Range("H19").Value = 3947.5
$ws.Range("I19").Value = 5000
$ws.Range("J19").Value = 2895
$ws.Range("K19").Value = 5000
$ws.Range("L19").Value = 2895
$ws.Range("M19").Value = -4830
$ws.Range("N19").Value = -3235

$ws.Range("H25").Value = 4903.2
$ws.Range("I25").Value = 6000
$ws.Range("J25").Value = 4629
$ws.Range("K25").Value = 6000
$ws.Range("L25").Value = 4629
$ws.Range("M25").Value = -5770
$ws.Range("N25").Value = -5089

$ws.Range("H28").Value = 5000
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 5000
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 5000
$ws.Range("N28").Value = -5464

$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 15000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 15000
$ws.Range("N31").Value = -15496

$ws.Range("H37").Value = 5000
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 5000
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 5000
$ws.Range("N37").Value = -5214

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 10000
$ws.Range("I13").Value = 10000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -9860

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()

$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 25000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630

$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 25000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184

$ws.Range("H81").Value = 1880
$ws.Range("I81").Value = 2050
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 4100
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -3039
$ws.Range("N81").Value = -4522

$ws.Range("H84").Value = 1880
$ws.Range("I84").Value = 2050
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 20500
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -15196
$ws.Range("N84").Value = -22608

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 1556
$ws.Range("I122").Value = 1832.3334
$ws.Range("J122").Value = 1390.2
$ws.Range("K122").Value = 5497.0002
$ws.Range("L122").Value = 4170.6
$ws.Range("M122").Value = -3047.0002
$ws.Range("N122").Value = -9070.6

$ws.Range("H141").Value = 87247.25
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 87247.25
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 87247.25
$ws.Range("N141").Value = -97607.25
